# Updated test data, JDE vendor search login configured in XML file.
#
# - The PhoneBook "Team" value for Joe Dejesus changes from
#   "Internal Medicine (FIM)" to "Family Practice (FIM)".
# - The PhoneBook sheet becomes the active/selected sheet, with cell G4
#   selected (previously VendorSearch was active with F6 selected).

$wb = $excel.ActiveWorkbook

$wsPhoneBook = $wb.Worksheets.Item("PhoneBook")

# Update the Team cell value for the phone book entry.
$wsPhoneBook.Range("G3").Value = "Family Practice (FIM)"

# Make PhoneBook the active sheet and select G4 on it.
# (VendorSearch's own selection, F6, is left untouched.)
$wsPhoneBook.Activate()
$wsPhoneBook.Range("G4").Select()
